# Fragments.xlsx - add references for new fragments (frg_PriceInformation,
# CarSalesContract_NO, Requisition_NO, VehicleControl_NO) to the Fragments
# usage matrix on sheet "Blad1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Header row (row 5): new document columns L, M, N ---
$ws.Range("L5").Value = "CarSalesContract_NO"
$ws.Range("M5").Value = "Requisition_NO"
$ws.Range("N5").Value = "VehicleControl_NO"

# --- Mark usage with "x" for CarSalesContract_NO column (L) ---
$ws.Range("L6").Value  = "x"
$ws.Range("L7").Value  = "x"
$ws.Range("L8").Value  = "x"
$ws.Range("L9").Value  = "x"
$ws.Range("L11").Value = "x"
$ws.Range("L14").Value = "x"
$ws.Range("L17").Value = "x"
$ws.Range("L18").Value = "x"
$ws.Range("L19").Value = "x"
$ws.Range("L21").Value = "x"
$ws.Range("L22").Value = "x"
$ws.Range("L23").Value = "x"
$ws.Range("L25").Value = "x"
$ws.Range("L26").Value = "x"
$ws.Range("L29").Value = "x"

# --- Mark usage with "x" for Requisition_NO column (M) ---
$ws.Range("M9").Value  = "x"
$ws.Range("M11").Value = "x"
$ws.Range("M17").Value = "x"
$ws.Range("M18").Value = "x"
$ws.Range("M19").Value = "x"
$ws.Range("M22").Value = "x"
$ws.Range("M23").Value = "x"
$ws.Range("M25").Value = "x"
$ws.Range("M26").Value = "x"
$ws.Range("M27").Value = "x"
$ws.Range("M28").Value = "x"

# --- New fragment row: frg_PriceInformation (row 29), used by CarPriceEstimate_NO (B) ---
$ws.Range("A29").Value = "frg_PriceInformation"
$ws.Range("B29").Value = "x"

# --- View/selection changes ---
$ws.Range("A4").Activate()
$ws.Range("N6").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
